# Applies the numeric updates to column F ("参与人数"/attendance-like counts)
# across the "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) sheets.
# "本地生活" (sheet3) is unchanged.

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 301
$ws1.Range("F8").Value  = 40
$ws1.Range("F9").Value  = 503
$ws1.Range("F11").Value = 293
$ws1.Range("F14").Value = 228
$ws1.Range("F17").Value = 6576
$ws1.Range("F19").Value = 69
$ws1.Range("F21").Value = 7517
$ws1.Range("F26").Value = 1300
$ws1.Range("F29").Value = 20
$ws1.Range("F31").Value = 66
$ws1.Range("F32").Value = 204
$ws1.Range("F33").Value = 193
$ws1.Range("F34").Value = 1606
$ws1.Range("F40").Value = 1701

# --- 演出 sheet ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 48

# --- 全部类型 sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value  = 301
$ws4.Range("F10").Value = 40
$ws4.Range("F11").Value = 503
$ws4.Range("F14").Value = 293
$ws4.Range("F18").Value = 228
$ws4.Range("F21").Value = 6576
$ws4.Range("F23").Value = 69
$ws4.Range("F25").Value = 7517
$ws4.Range("F30").Value = 1348
$ws4.Range("F33").Value = 20
$ws4.Range("F35").Value = 66
$ws4.Range("F36").Value = 48
$ws4.Range("F37").Value = 204
$ws4.Range("F38").Value = 193
$ws4.Range("F39").Value = 1606
$ws4.Range("F45").Value = 1701
